$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Shift the six daily-date blocks forward by one week (2023-09-18..23 -> 2023-09-25..30)
$dates = @(
    "2023-09-25",
    "2023-09-26",
    "2023-09-27",
    "2023-09-28",
    "2023-09-29",
    "2023-09-30"
)

$startRows = @(2, 11, 20, 29, 38, 47)
$endRows   = @(10, 19, 28, 37, 46, 55)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r1 = $startRows[$i]
    $r2 = $endRows[$i]
    $ws.Range("A$r1`:A$r2").Value = $dates[$i]
}

# Update the saved active selection to B48
$ws.Range("B48").Select()
